$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column L, mirroring column K (year 2021) with the 2022 data.
# First copy each existing K-column cell's formatting down into L so the
# new cells pick up the right style, then set the new values.

# Row 2: empty bottom-border cell, no value, same style as K2.
$ws.Range("K2").Copy()
$ws.Range("L2").PasteSpecial(-4122)

# Row 3: header year value.
$ws.Range("K3").Copy()
$ws.Range("L3").PasteSpecial(-4122)
$ws.Range("L3").Value = 2022

# Row 4.
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("L4").Value = 370

# Row 5.
$ws.Range("K5").Copy()
$ws.Range("L5").PasteSpecial(-4122)
$ws.Range("L5").Value = 137

# Row 6.
$ws.Range("K6").Copy()
$ws.Range("L6").PasteSpecial(-4122)
$ws.Range("L6").Value = 314

# Row 7.
$ws.Range("K7").Copy()
$ws.Range("L7").PasteSpecial(-4122)
$ws.Range("L7").Value = 121

# Row 8: this one gets a thousands-separator number format (new style).
$ws.Range("K8").Copy()
$ws.Range("L8").PasteSpecial(-4122)
$ws.Range("L8").NumberFormat = "#,##0"
$ws.Range("L8").Value = 50

# Row 9.
$ws.Range("K9").Copy()
$ws.Range("L9").PasteSpecial(-4122)
$ws.Range("L9").Value = 16

# Move the selection to L2, matching the author's final cursor position.
$ws.Range("L2").Select()
